$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as plain text so values
# like "1.00", "298.67" or "0.0779" are not coerced into numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.211.47'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.226.74'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -1.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '298.67'
$ws.Range('E5').Value = '  -2.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '90.62'
$ws.Range('E6').Value = '  -4.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.558'
$ws.Range('E7').Value = '  -2.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -5.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.23'
$ws.Range('E10').Value = '  -4.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0779'
$ws.Range('E11').Value = '  -3.00%  '
$ws.Range('E12').Value = '  -3.32%  '
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.565.85'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.226.06'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.40'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.778'
$ws.Range('E17').Value = '  -6.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.068.27'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.34'
$ws.Range('E19').Value = '  +3.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0908'
$ws.Range('E20').Value = '  -4.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.00'
$ws.Range('E21').Value = '  -5.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.33'
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.49'
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('E24').Value = '  -4.85%  '
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('E26').Value = '  -6.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.39'
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.27'
$ws.Range('E30').Value = '  -3.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '151.13'
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('E32').Value = '  -8.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0764'
$ws.Range('E33').Value = '  -4.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.50'
$ws.Range('E34').Value = '  -5.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.116'
$ws.Range('E35').Value = '  -2.08%  '
$ws.Range('E36').Value = '  -5.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.85'
$ws.Range('E37').Value = '  -7.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.68'
$ws.Range('E38').Value = '  -7.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0300'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('E40').Value = '  -6.95%  '
$ws.Range('E41').Value = '  -4.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.47'
$ws.Range('E42').Value = '  -9.75%  '
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.796.62'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.82'
$ws.Range('E45').Value = '  +9.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.185'
$ws.Range('E46').Value = '  -3.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '67.99'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '94.77'
$ws.Range('E48').Value = '  -4.45%  '
$ws.Range('E49').Value = '  -7.55%  '
$ws.Range('E50').Value = '  -4.59%  '
$ws.Range('E51').Value = '  -5.96%  '

# Row 27/28 swap: InjectiveProtocol moves to row 27, Toncoin moves to row 28
$ws.Range('B27').Value = 'InjectiveProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '39.03'
$ws.Range('E27').Value = '  +2.41%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('E28').Value = '  -0.86%  '
